$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the product/brand labels (shared strings) ---
# Column headers (row 1) and matching row labels (column A) must stay in sync.

$ws.Range("C1").Value = "ดัชมิลล์"
$ws.Range("A3").Value = "ดัชมิลล์"

$ws.Range("D1").Value = "โฟร์โมสต์"
$ws.Range("A4").Value = "โฟร์โมสต์"

$ws.Range("E1").Value = "แดรี่โฮม"
$ws.Range("A5").Value = "แดรี่โฮม"

$ws.Range("G1").Value = "เอ็มมิลค์"
$ws.Range("A7").Value = "เอ็มมิลค์"

$ws.Range("H1").Value = "ไทยเดนมาร์ค"
$ws.Range("A8").Value = "ไทยเดนมาร์ค"

# --- Update the Jaccard similarity / cross-validation matrix values ---

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.1108032196764444
$ws.Range("D2").Value = 0.09446707970288731
$ws.Range("E2").Value = 0.01956254241693767
$ws.Range("F2").Value = 0.02964585807536208
$ws.Range("G2").Value = 0.04017087007589609
$ws.Range("H2").Value = 0.05003973575996431
$ws.Range("B3").Value = 0.1108032196764444
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.124671961322111
$ws.Range("E3").Value = 0.02189330424219897
$ws.Range("F3").Value = 0.04173466467630537
$ws.Range("G3").Value = 0.0004553734061930783
$ws.Range("H3").Value = 0.1114107448998216
$ws.Range("B4").Value = 0.09446707970288731
$ws.Range("C4").Value = 0.124671961322111
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.03547041304011841
$ws.Range("G4").Value = 0.0002579979360165118
$ws.Range("H4").Value = 0.117420824480639
$ws.Range("B5").Value = 0.01956254241693767
$ws.Range("C5").Value = 0.02189330424219897
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.01346409356032294
$ws.Range("G5").Value = 0.003787878787878788
$ws.Range("H5").Value = 0.0426602026144021
$ws.Range("B6").Value = 0.02964585807536208
$ws.Range("C6").Value = 0.04173466467630537
$ws.Range("D6").Value = 0.03547041304011841
$ws.Range("E6").Value = 0.01346409356032294
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.02487639490014739
$ws.Range("B7").Value = 0.04017087007589609
$ws.Range("C7").Value = 0.0004553734061930783
$ws.Range("D7").Value = 0.0002579979360165118
$ws.Range("E7").Value = 0.003787878787878788
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.00004782629489693433
$ws.Range("B8").Value = 0.05003973575996431
$ws.Range("C8").Value = 0.1114107448998216
$ws.Range("D8").Value = 0.117420824480639
$ws.Range("E8").Value = 0.0426602026144021
$ws.Range("F8").Value = 0.02487639490014739
$ws.Range("G8").Value = 0.00004782629489693433
$ws.Range("H8").Value = 1
